# Daily attendance processing - reorders the "Recorded By" (column G)
# value on every row: when it starts with "System," (System recorded the
# session first, followed by the human reviewer(s)), flip the order so the
# human reviewer(s) are listed first and "System" last.
#
# Example: "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#          "System, system, backup@backdoor.com" -> "backup@backdoor.com, system, System"
#
# Rows whose "Recorded By" value does not start with "System," are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ($text.StartsWith("System,")) {
        $parts = $text.Split(",")
        $count = $parts.Count

        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i].Trim()
        }

        $cell.Value2 = [string]::Join(", ", $reversed)
    }
}
